$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Prep rows 7 & 8: clone row 6's formatting (borders/wrap/vertical-align)
# down into the two new BOM rows before we populate them. Using
# xlPasteFormats keeps the paste re-using the workbook's existing style
# records instead of minting new ones.
$ws.Range("A6:F6").Copy()
$ws.Range("A7:F8").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Row 2: decoupling caps now cover C1-C12 (was C1, C2) -----------------
$ws.Range("A2").Value = 12
$ws.Range("D2").Value = "'C1, C2, C3, C4, C5, C6, C7, C8, C9, C10, C11, C12"
$ws.Range("F2").Value = 12
$ws.Rows.Item(2).RowHeight = 45

# --- Row 6: was AM26C32IPWR/U2, now the new 5V->3.3V level translator -----
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "'SN74LVC1T45DBVR"
$ws.Range("C6").Value = "'296-16843-1-ND"
$ws.Range("D6").Value = "'U2, U3, U5, U6"
$ws.Range("F6").Value = 4

# --- Row 7: the old AM26C32IPWR/U2 part, redesignated to U4 ---------------
$ws.Range("A7").Value = 1
$ws.Range("B7").Value = "'AM26C32IPWR"
$ws.Range("C7").Value = "'296-30088-1-ND"
$ws.Range("D7").Value = "'U4"
$ws.Range("F7").Value = 1

# --- Row 8: new 3.3V regulator (U7) ----------------------------------------
$ws.Range("A8").Value = 1
$ws.Range("B8").Value = "'TLV73333PDBVR"
$ws.Range("C8").Value = "'296-40673-1-ND"
$ws.Range("D8").Value = "'U7"
$ws.Range("F8").Value = 1

$ws.Range("E11").Select()
